{"js": "const replacements = [\n  [\"2025-03-19 Wednesday\", \"2025-03-20 Thursday\"],\n  [\"477\u00d79=\", \"909\u00d74=\"],\n  [\"472\u00d77=\", \"626\u00d75=\"],\n  [\"591\u00d77=\", \"690\u00d72=\"],\n  [\"987\u00d78=\", \"435\u00d72=\"],\n  [\"652\u00d74=\", \"392\u00d77=\"],\n  [\"511\u00d77=\", \"482\u00d72=\"],\n  [\"354\u00d77=\", \"825\u00d77=\"],\n  [\"252\u00d77=\", \"151\u00d73=\"],\n  [\"959\u00d73=\", \"832\u00d76=\"],\n  [\"148\u00d78=\", \"563\u00d76=\"],\n  [\"641\u00d78=\", \"726\u00d74=\"],\n  [\"232\u00d78=\", \"470\u00d77=\"],\n  [\"322\u00d74=\", \"674\u00d73=\"],\n  [\"689\u00d72=\", \"441\u00d72=\"],\n  [\"498\u00d79=\", \"341\u00d79=\"],\n  [\"889\u00d75=\", \"931\u00d79=\"],\n  [\"212\u00d76=\", \"358\u00d78=\"],\n  [\"626\u00d77=\", \"797\u00d79=\"],\n  [\"222\u00d74=\", \"743\u00d79=\"],\n  [\"114\u00d77=\", \"217\u00d76=\"],\n  [\"820\u00d73=\", \"427\u00d75=\"],\n  [\"596\u00d76=\", \"420\u00d77=\"],\n  [\"694\u00d76=\", \"247\u00d78=\"],\n  [\"637\u00d72=\", \"810\u00d75=\"],\n];\n\nfor (const [searchText, replaceText] of replacements) {\n  const results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-19 Wednesday\", \"2025-03-20 Thursday\"),\n    @(\"477\u00d79=\", \"909\u00d74=\"),\n    @(\"472\u00d77=\", \"626\u00d75=\"),\n    @(\"591\u00d77=\", \"690\u00d72=\"),\n    @(\"987\u00d78=\", \"435\u00d72=\"),\n    @(\"652\u00d74=\", \"392\u00d77=\"),\n    @(\"511\u00d77=\", \"482\u00d72=\"),\n    @(\"354\u00d77=\", \"825\u00d77=\"),\n    @(\"252\u00d77=\", \"151\u00d73=\"),\n    @(\"959\u00d73=\", \"832\u00d76=\"),\n    @(\"148\u00d78=\", \"563\u00d76=\"),\n    @(\"641\u00d78=\", \"726\u00d74=\"),\n    @(\"232\u00d78=\", \"470\u00d77=\"),\n    @(\"322\u00d74=\", \"674\u00d73=\"),\n    @(\"689\u00d72=\", \"441\u00d72=\"),\n    @(\"498\u00d79=\", \"341\u00d79=\"),\n    @(\"889\u00d75=\", \"931\u00d79=\"),\n    @(\"212\u00d76=\", \"358\u00d78=\"),\n    @(\"626\u00d77=\", \"797\u00d79=\"),\n    @(\"222\u00d74=\", \"743\u00d79=\"),\n    @(\"114\u00d77=\", \"217\u00d76=\"),\n    @(\"820\u00d73=\", \"427\u00d75=\"),\n    @(\"596\u00d76=\", \"420\u00d77=\"),\n    @(\"694\u00d76=\", \"247\u00d78=\"),\n    @(\"637\u00d72=\", \"810\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
